$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 90
$ws.Range("C5").Value = 0.6325197058534924
$ws.Range("D5").Value = 0.02765385105211411
$ws.Range("E5").Value = 0.8661309802792049
